$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30; this shifts the existing rows 30-62
# down to 31-63 (matching the dimension growing from A1:R62 to A1:R63).
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly price record.
$ws.Cells.Item(30, 1).Value = 8
$ws.Cells.Item(30, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(30, 3).Value = "Coquimbo"
$ws.Cells.Item(30, 4).Value = 44810
$ws.Cells.Item(30, 5).Value = 4
$ws.Cells.Item(30, 6).Value = 100114007
$ws.Cells.Item(30, 7).Value = "Jengibre"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 540
$ws.Cells.Item(30, 11).Value = 14000
$ws.Cells.Item(30, 12).Value = 15000
$ws.Cells.Item(30, 13).Value = 14500
$ws.Cells.Item(30, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(30, 15).Value = "Perú"
$ws.Cells.Item(30, 16).Value = 1115
$ws.Cells.Item(30, 17).Value = 13
$ws.Cells.Item(30, 18).Value = "Hortaliza"
